$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range("D2").Value = '30.541.42'
$ws.Range("E2").Value = '  +0.42%  '
$ws.Range("D3").Value = '1.871.95'
$ws.Range("E3").Value = '  +0.31%  '
$ws.Range("D4").Value = '1.000'
$ws.Range("E4").Value = '  +0.02%  '
$ws.Range("D5").Value = '235.52'
$ws.Range("E5").Value = '  -0.56%  '
$ws.Range("D6").Value = '1.000'
$ws.Range("E6").Value = '  +0.03%  '
$ws.Range("D7").Value = '0.4829'
$ws.Range("E7").Value = '  +0.50%  '
$ws.Range("D8").Value = '0.2809'
$ws.Range("E8").Value = '  +0.07%  '
$ws.Range("D9").Value = '0.06515'
$ws.Range("E9").Value = '  -0.12%  '
$ws.Range("D10").Value = '1.903.02'
$ws.Range("E10").Value = '  +2.01%  '
$ws.Range("D11").Value = '0.07446'
$ws.Range("E11").Value = '  +0.08%  '
$ws.Range("D12").Value = '16.35'
$ws.Range("E12").Value = '  +0.35%  '
$ws.Range("D13").Value = '5.082'
$ws.Range("E13").Value = '  +0.20%  '
$ws.Range("D14").Value = '87.32'
$ws.Range("E14").Value = '  -1.05%  '
$ws.Range("D15").Value = '0.6483'
$ws.Range("E15").Value = '  -1.00%  '
$ws.Range("D16").Value = '30.506.92'
$ws.Range("E16").Value = '  +0.41%  '
$ws.Range("D17").Value = '1.000'
$ws.Range("E17").Value = '  -0.02%  '
$ws.Range("D18").Value = '13.02'
$ws.Range("E18").Value = '  -1.96%  '
$ws.Range("B19").Value = 'ShibaInu'
$ws.Range("C19").Value = 'https://coinranking.com/coin/xz24e0BjL+shibainu-shib'
$ws.Range("D19").Value = '0.000007550'
$ws.Range("E19").Value = '  -0.83%  '
$ws.Range("B20").Value = 'BitcoinCash'
$ws.Range("C20").Value = 'https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch'
$ws.Range("D20").Value = '231.77'
$ws.Range("E20").Value = '  +5.67%  '
$ws.Range("B21").Value = 'WrappedliquidstakedEther2.0'
$ws.Range("C21").Value = 'https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth'
$ws.Range("D21").Value = '2.104.87'
$ws.Range("E21").Value = '  -0.23%  '
$ws.Range("D22").Value = '1.000'
$ws.Range("E22").Value = '  +0.03%  '
$ws.Range("D23").Value = '5.168'
$ws.Range("E23").Value = '  -2.06%  '
$ws.Range("D24").Value = '6.120'
$ws.Range("E24").Value = '  -0.58%  '
$ws.Range("D25").Value = '9.360'
$ws.Range("E25").Value = '  +1.06%  '
$ws.Range("D26").Value = '167.37'
$ws.Range("E26").Value = '  +1.59%  '
$ws.Range("D27").Value = '18.40'
$ws.Range("E27").Value = '  -0.89%  '
$ws.Range("D28").Value = '1.924'
$ws.Range("E28").Value = '  -1.10%  '
$ws.Range("D29").Value = '0.1040'
$ws.Range("E29").Value = '  +11.76%  '
$ws.Range("D30").Value = '1.377'
$ws.Range("E30").Value = '  -4.91%  '
$ws.Range("D31").Value = '4.275'
$ws.Range("E31").Value = '  -0.30%  '
$ws.Range("D32").Value = '4.000'
$ws.Range("E32").Value = '  -0.09%  '
$ws.Range("D33").Value = '0.04977'
$ws.Range("E33").Value = '  -1.33%  '
$ws.Range("D34").Value = '1.185'
$ws.Range("E34").Value = '  -0.08%  '
$ws.Range("D35").Value = '0.7378'
$ws.Range("E35").Value = '  -1.66%  '
$ws.Range("D36").Value = '0.9997'
$ws.Range("E36").Value = '  +0.17%  '
$ws.Range("D37").Value = '2.714'
$ws.Range("E37").Value = '  +0.64%  '
$ws.Range("D38").Value = '0.01919'
$ws.Range("E38").Value = '  +5.02%  '
$ws.Range("D39").Value = '2.639'
$ws.Range("E39").Value = '  +0.73%  '
$ws.Range("D40").Value = '0.9136'
$ws.Range("E40").Value = '  +1.37%  '
$ws.Range("D41").Value = '2.051'
$ws.Range("E41").Value = '  -0.86%  '
$ws.Range("D42").Value = '106.08'
$ws.Range("E42").Value = '  -0.42%  '
$ws.Range("B43").Value = 'PaxDollar'
$ws.Range("C43").Value = 'https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp'
$ws.Range("D43").Value = '0.9974'
$ws.Range("E43").Value = '  -0.42%  '
$ws.Range("B44").Value = 'TheSandbox'
$ws.Range("C44").Value = 'https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand'
$ws.Range("D44").Value = '0.4221'
$ws.Range("E44").Value = '  -1.48%  '
$ws.Range("D45").Value = '5.611'
$ws.Range("E45").Value = '  -5.36%  '
$ws.Range("D46").Value = '7.289'
$ws.Range("E46").Value = '  -1.56%  '
$ws.Range("D47").Value = '62.99'
$ws.Range("E47").Value = '  -2.80%  '
$ws.Range("D48").Value = '0.1247'
$ws.Range("E48").Value = '  -3.56%  '
$ws.Range("D49").Value = '8.926'
$ws.Range("E49").Value = '  -0.39%  '
$ws.Range("D50").Value = '1.445'
$ws.Range("E50").Value = '  -1.99%  '
$ws.Range("D51").Value = '33.74'
$ws.Range("E51").Value = '  -1.45%  '
